$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (1008, f, 8700, PRJ-03) below the existing table
$ws.Range("A11").Value = 1008
$ws.Range("B11").Value = "f"
$ws.Range("C11").Value = 8700
$ws.Range("D11").Value = "PRJ-03"

# Move the active selection to A12, as if the user pressed Enter after
# typing the new row
$ws.Range("A12").Select()
